$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-12 17:15:02", 0.0008),
    @("2023-12-12 17:15:17", 0.001),
    @("2023-12-12 17:15:23", 0.0004)
)

$startRow = 247
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
